$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: add the date and the event description that were previously blank
$ws.Range("A7").Value = (Get-Date -Year 2020 -Month 4 -Day 2 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("B7").Value = "Début de la rédaction de la documentation du projet"

# Update the selected cell to match the new active selection
$ws.Range("A14").Select()
